# Reposition the six "card" groups (icon + text box) on slide 4
# from a single-row layout into a 2-row x 3-column grid.
#
# Shape.Left/.Top are expressed in points; PowerPoint stores the
# underlying OOXML <a:off> in EMU (1 pt = 12700 EMU) using Single
# (32-bit float) precision, so the literals below are chosen to land
# on the exact target EMU value after that round-trip.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Group 13: off 2976003,1535122 -> 4690126,1665810 EMU
$shp = $s.Shapes.Item("Group 13")
$shp.Left = 369.30126953125
$shp.Top = 131.16615295410156

# Group 18: off 12206614,1535122 -> 4690126,5314277 EMU
$shp = $s.Shapes.Item("Group 18")
$shp.Left = 369.30126953125
$shp.Top = 418.447021484375

# Group 23: off 6128505,5186060 -> 13830491,5346764 EMU
$shp = $s.Shapes.Item("Group 23")
$shp.Left = 1089.0150146484375
$shp.Top = 421.00506591796875

# Group 28: off 9327360,1535122 -> 13953560,1775472 EMU
$shp = $s.Shapes.Item("Group 28")
$shp.Left = 1098.70556640625
$shp.Top = 139.8009490966797

# Group 33: off 15273652,1535122 -> 8748010,5419405 EMU
$shp = $s.Shapes.Item("Group 33")
$shp.Left = 688.8197021484375
$shp.Top = 426.7248229980469

# Group 38: off 6257913,1535122 -> 8681075,1665810 EMU
$shp = $s.Shapes.Item("Group 38")
$shp.Left = 683.5492553710938
$shp.Top = 131.16615295410156

# TextBox 48: off 2976003,2641282 -> 4690126,2771971 EMU
$shp = $s.Shapes.Item("TextBox 48")
$shp.Left = 369.30126953125
$shp.Top = 218.26544189453125

# TextBox 49: off 6257913,2641282 -> 8681075,2771971 EMU
$shp = $s.Shapes.Item("TextBox 49")
$shp.Left = 683.5492553710938
$shp.Top = 218.26544189453125

# TextBox 50: off 9204290,2617827 -> 13830491,2858177 EMU
$shp = $s.Shapes.Item("TextBox 50")
$shp.Left = 1089.0150146484375
$shp.Top = 225.05331420898438

# TextBox 51: off 12206614,2641282 -> 4690126,6420437 EMU
$shp = $s.Shapes.Item("TextBox 51")
$shp.Left = 369.30126953125
$shp.Top = 505.5462341308594

# TextBox 52: off 15206717,2641282 -> 8681075,6525566 EMU
$shp = $s.Shapes.Item("TextBox 52")
$shp.Left = 683.5492553710938
$shp.Top = 513.8240966796875

# TextBox 54: off 6128505,6259733 -> 13830491,6420437 EMU
$shp = $s.Shapes.Item("TextBox 54")
$shp.Left = 1089.0150146484375
$shp.Top = 505.5462341308594
